$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# from 45175 (2023-09-06) to 45177 (2023-09-08) for every data row
# (rows 2 through 358).
$ws.Range("C2:C358").Value = 45177
